$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new rows (394-415) of player stats, matching the source workbook upload
$newRows = @(
    @('Philipe', 1, 3, 2, 3, 1, 0, 1, 0, 0),
    @('Jorge', 1, 3, 2, 0, 1, 0, 1, 0, 0),
    @('Michel', 1, 3, 2, 3, 1, 0, 1, 0, 0),
    @('Marcelão', 1, 3, 2, 1, 1, 0, 1, 0, 0),
    @('Juscielio', 1, 3, 2, 1, 1, 0, 1, 0, 0),
    @('Euler', 3, 2, 2, 0, 1, 1, 0, 0, 0),
    @('Athos', 3, 2, 2, 2, 1, 1, 0, 0, 0),
    @('David', 3, 2, 2, 2, 1, 1, 0, 0, 0),
    @('Corinthiano', 3, 2, 2, 1, 1, 1, 0, 0, 0),
    @('Leandrinho', 3, 2, 2, 2, 1, 1, 0, 0, 0),
    @('Guinha', 3, 2, 2, 0, 1, 0, 0, 0, 0),
    @('Peixe', 3, 2, 2, 1, 1, 0, 0, 0, 0),
    @('Ismael', 3, 2, 2, 0, 1, 0, 0, 0, 0),
    @('Eder', 3, 2, 2, 2, 1, 0, 0, 0, 0),
    @('Boneco', 3, 2, 2, 1, 1, 0, 0, 0, 0),
    @('Fernando', 2, 1, 3, 2, 1, 0, 0, 0, 0),
    @('Romario', 2, 1, 3, 3, 1, 0, 0, 0, 0),
    @('Du', 2, 1, 3, 0, 1, 0, 0, 0, 0),
    @('Leandrão', 2, 1, 3, 0, 1, 0, 0, 0, 0),
    @('Cabeleira', 2, 1, 3, 1, 1, 0, 0, 0, 0),
    @('Matheus', 4, 4, 5, 0, 1, 0, 1, 1, 11),
    @('Igor Goleiro', 5, 4, 4, 0, 1, 1, 0, 0, 12),
)

$startRow = 394
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]   # A: Jogadores (name)
    $ws.Cells.Item($r, 3).Value = $data[1]   # C: Vitorias
    $ws.Cells.Item($r, 4).Value = $data[2]   # D: Empate
    $ws.Cells.Item($r, 5).Value = $data[3]   # E: Derrotas
    $ws.Cells.Item($r, 6).Value = $data[4]   # F: Gols
    $ws.Cells.Item($r, 7).Value = $data[5]   # G: Partidas
    $ws.Cells.Item($r, 8).Value = $data[6]   # H: Tarde de Vitoria
    $ws.Cells.Item($r, 9).Value = $data[7]   # I: La barca
    $ws.Cells.Item($r, 10).Value = $data[8]  # J: Craque do Dia
    $ws.Cells.Item($r, 11).Value = $data[9]  # K: Gols Sofridos
}

# Update the selected cell / view to reflect appended data (matches saved workbook state)
$ws.Range("A416").Select() | Out-Null
